# Hjemme passive tweaks lichtwark deleted values
#
# Columns B:E (Subj 1-4, header row 1 = 1,2,3,4) are overwritten with the
# data that otherwise lives under the "16" / "20" header columns (O/R on the
# CON row, O/R on the STR row), i.e. the original "deleted" subjects are
# replaced by re-pasted values for subjects 16/20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - subject/header ids
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 - CON
$ws.Range("B2").Value = 10.466253638017726
$ws.Range("C2").Value = 5.7490026772609859
$ws.Range("D2").Value = 10.663079852511585
$ws.Range("E2").Value = 8.9113260483376351

# Row 3 - STR
$ws.Range("B3").Value = 6.0993965164398682
$ws.Range("C3").Value = 8.0333669674263248
$ws.Range("D3").Value = 7.2841762501876959
$ws.Range("E3").Value = 8.686184479998504

# Narrow the active selection from the whole used range down to B1:E3
$ws.Range("B1:E3").Select
